$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending / Target clusters, in order
$clusters = @("ECs", "FAPs", "M2", "sCs")

# Per-sending-cluster stats: G,H,I,J (Ligand stats) - same for all 4 target rows of that sending cluster
$ligandStats = @{
    "ECs"  = @(7.079689000000001, 21.239067, 0.5033576067109902, 0.5033576067109902)
    "FAPs" = @(4.058683666666667, 12.176051, 0.2885676612136944, 0.2885676612136945)
    "M2"   = @(0.463315, 1.389945, 0.03294115455541936, 0.03294115455541936)
    "sCs"  = @(2.463241333333333, 7.389724, 0.175133577519896, 0.175133577519896)
}

# Per-target-cluster stats: K,L,M,N,O,P (Receptor stats) - same regardless of sending cluster
$receptorStats = @{
    "ECs"  = @(2, 0.6666666666666666, 23.34377233333333, 70.031317, 0.6243713731385436, 0.6243713731385436)
    "FAPs" = @(3, 1, 8.314670666666666, 24.944012, 0.2223908915496236, 0.2223908915496236)
    "M2"   = @(1, 0.3333333333333333, 0.02865466666666667, 0.085964, 0.0007664208388438813, 0.0007664208388438813)
    "sCs"  = @(3, 1, 5.700542666666667, 17.101628, 0.1524713144729888, 0.1524713144729888)
}

# Per (sending,target) pair: Q,R,S,T (edge stats)
$edgeStats = @{
    "ECs|ECs"   = @(165.2666482068043, 1487.399833861239, 0.3142820800818719, 0.3142820800818719)
    "ECs|FAPs"  = @(58.86528245742267, 529.787542116804, 0.1119421469247419, 0.1119421469247419)
    "ECs|M2"    = @(0.2028661283986667, 1.825795155588, 0.0003857837591738855, 0.0003857837591738855)
    "ECs|sCs"   = @(40.35806921123068, 363.2226229010761, 0.0767475959452024, 0.0767475959452024)
    "FAPs|ECs"  = @(94.74498748768524, 852.7048873891671, 0.1801733868753725, 0.1801733868753725)
    "FAPs|FAPs" = @(33.74661802851244, 303.7195622566121, 0.06417481944970324, 0.06417481944970324)
    "FAPs|M2"   = @(0.1163002275737778, 1.046702048164, 0.0002211642689706166, 0.0002211642689706167)
    "FAPs|sCs"  = @(23.13669941233645, 208.2302947110281, 0.0439982906196481, 0.04399829061964811)
    "M2|ECs"    = @(10.81551987861833, 97.339678907565, 0.02056751390253618, 0.02056751390253618)
    "M2|FAPs"   = @(3.852311639926666, 34.67080475934, 0.007325812730253656, 0.007325812730253656)
    "M2|M2"     = @(0.01327613688666666, 0.11948523198, 0.00002524678730685045, 0.00002524678730685045)
    "M2|sCs"    = @(2.641146925606667, 23.77032233046, 0.005022581135322674, 0.005022581135322674)
    "sCs|ECs"   = @(57.50134488738978, 517.512103986508, 0.109348392278763, 0.109348392278763)
    "sCs|FAPs"  = @(20.48104045918755, 184.329364132688, 0.03894811244492478, 0.03894811244492478)
    "sCs|M2"    = @(0.07058335932622221, 0.635250233936, 0.0001342260233925286, 0.0001342260233925286)
    "sCs|sCs"   = @(14.04181231896356, 126.376310870672, 0.02670284677281562, 0.02670284677281562)
}

$rowNum = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ls = $ligandStats[$sending]
        $rs = $receptorStats[$target]
        $es = $edgeStats["$sending|$target"]

        $ws.Cells.Item($rowNum, 1).Value = $sending
        $ws.Cells.Item($rowNum, 2).Value = "Efnb1"
        $ws.Cells.Item($rowNum, 3).Value = "Ephb4"
        $ws.Cells.Item($rowNum, 4).Value = $target
        $ws.Cells.Item($rowNum, 5).Value = 3
        $ws.Cells.Item($rowNum, 6).Value = 1
        $ws.Cells.Item($rowNum, 7).Value = $ls[0]
        $ws.Cells.Item($rowNum, 8).Value = $ls[1]
        $ws.Cells.Item($rowNum, 9).Value = $ls[2]
        $ws.Cells.Item($rowNum, 10).Value = $ls[3]
        $ws.Cells.Item($rowNum, 11).Value = $rs[0]
        $ws.Cells.Item($rowNum, 12).Value = $rs[1]
        $ws.Cells.Item($rowNum, 13).Value = $rs[2]
        $ws.Cells.Item($rowNum, 14).Value = $rs[3]
        $ws.Cells.Item($rowNum, 15).Value = $rs[4]
        $ws.Cells.Item($rowNum, 16).Value = $rs[5]
        $ws.Cells.Item($rowNum, 17).Value = $es[0]
        $ws.Cells.Item($rowNum, 18).Value = $es[1]
        $ws.Cells.Item($rowNum, 19).Value = $es[2]
        $ws.Cells.Item($rowNum, 20).Value = $es[3]

        $rowNum++
    }
}
